{"js": "// Replace the 25 division-problem cell texts with their updated values.\n// Pairs are listed in document order (top-to-bottom, left-to-right within a\n// row) so that each `search()` targets the original, not-yet-modified text\n// even though a couple of the new values coincide with old values used\n// elsewhere in the table (e.g. \"115\u00f72=\" is both an old value at the top of\n// the table and a new value further down).\nconst replacements = [\n  [\"646\u00f77=\", \"689\u00f74=\"],\n  [\"442\u00f76=\", \"446\u00f79=\"],\n  [\"115\u00f72=\", \"525\u00f77=\"],\n  [\"136\u00f77=\", \"575\u00f76=\"],\n  [\"798\u00f74=\", \"995\u00f74=\"],\n  [\"775\u00f76=\", \"744\u00f74=\"],\n  [\"190\u00f77=\", \"226\u00f79=\"],\n  [\"620\u00f72=\", \"157\u00f72=\"],\n  [\"351\u00f78=\", \"761\u00f76=\"],\n  [\"730\u00f76=\", \"469\u00f78=\"],\n  [\"611\u00f74=\", \"472\u00f76=\"],\n  [\"789\u00f74=\", \"688\u00f75=\"],\n  [\"215\u00f78=\", \"884\u00f73=\"],\n  [\"956\u00f74=\", \"910\u00f74=\"],\n  [\"688\u00f79=\", \"619\u00f73=\"],\n  [\"638\u00f73=\", \"115\u00f72=\"],\n  [\"757\u00f76=\", \"843\u00f72=\"],\n  [\"983\u00f72=\", \"124\u00f75=\"],\n  [\"705\u00f78=\", \"153\u00f72=\"],\n  [\"636\u00f72=\", \"115\u00f75=\"],\n  [\"148\u00f75=\", \"893\u00f76=\"],\n  [\"432\u00f78=\", \"975\u00f73=\"],\n  [\"376\u00f72=\", \"705\u00f79=\"],\n  [\"602\u00f78=\", \"379\u00f73=\"],\n  [\"389\u00f76=\", \"515\u00f77=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace only the first (and expected-only) match for this exact text.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the 25 division-problem cell texts with their updated values.\n# Pairs are listed in document order (top-to-bottom, left-to-right within a\n# row) so that each Find targets the original, not-yet-modified text even\n# though a couple of the new values coincide with old values used elsewhere\n# in the table (e.g. \"115\u00f72=\" is both an old value at the top of the\n# table and a new value further down). Each old value occurs exactly once\n# in the original document, so replacing one-at-a-time in this order is\n# unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"646\u00f77=\"; New = \"689\u00f74=\" },\n    @{ Old = \"442\u00f76=\"; New = \"446\u00f79=\" },\n    @{ Old = \"115\u00f72=\"; New = \"525\u00f77=\" },\n    @{ Old = \"136\u00f77=\"; New = \"575\u00f76=\" },\n    @{ Old = \"798\u00f74=\"; New = \"995\u00f74=\" },\n    @{ Old = \"775\u00f76=\"; New = \"744\u00f74=\" },\n    @{ Old = \"190\u00f77=\"; New = \"226\u00f79=\" },\n    @{ Old = \"620\u00f72=\"; New = \"157\u00f72=\" },\n    @{ Old = \"351\u00f78=\"; New = \"761\u00f76=\" },\n    @{ Old = \"730\u00f76=\"; New = \"469\u00f78=\" },\n    @{ Old = \"611\u00f74=\"; New = \"472\u00f76=\" },\n    @{ Old = \"789\u00f74=\"; New = \"688\u00f75=\" },\n    @{ Old = \"215\u00f78=\"; New = \"884\u00f73=\" },\n    @{ Old = \"956\u00f74=\"; New = \"910\u00f74=\" },\n    @{ Old = \"688\u00f79=\"; New = \"619\u00f73=\" },\n    @{ Old = \"638\u00f73=\"; New = \"115\u00f72=\" },\n    @{ Old = \"757\u00f76=\"; New = \"843\u00f72=\" },\n    @{ Old = \"983\u00f72=\"; New = \"124\u00f75=\" },\n    @{ Old = \"705\u00f78=\"; New = \"153\u00f72=\" },\n    @{ Old = \"636\u00f72=\"; New = \"115\u00f75=\" },\n    @{ Old = \"148\u00f75=\"; New = \"893\u00f76=\" },\n    @{ Old = \"432\u00f78=\"; New = \"975\u00f73=\" },\n    @{ Old = \"376\u00f72=\"; New = \"705\u00f79=\" },\n    @{ Old = \"602\u00f78=\"; New = \"379\u00f73=\" },\n    @{ Old = \"389\u00f76=\"; New = \"515\u00f77=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $pair.Old,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceOne)\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
